$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Willette Pardie): update email and joining_date
$ws.Range("C4").Value = "infinitelooprogramming@gmail.com"
$ws.Range("E4").Value = 45531

# Row 3 (Cobby Jackett): update email and dob; convert to a mailto hyperlink
$ws.Range("C3").Value = "patoliyabhi17@gmail.com"
$ws.Range("C3").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:patoliyabhi17@gmail.com") | Out-Null
$ws.Range("C3").Style = "Hyperlink"
$ws.Range("D3").Value = 45531

# Update the last-active selection to match the final state
$ws.Range("D4").Select() | Out-Null
